$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New daily case/death records for Sergipe, rows 36-48 (01/05/2022 - 13/05/2022)
# Columns: A=Date(serial), B=col2, C=col3, D=col4, E=col5, F=col6
$data = @(
    @(44682, 0, 327081, 6343, 5, 0),
    @(44683, 0, 327089, 6343, 8, 0),
    @(44684, 0, 327099, 6344, 10, 1),
    @(44685, 0, 327106, 6344, 7, 0),
    @(44686, 0, 327114, 6344, 13, 0),
    @(44687, 0, 327132, 6344, 18, 0),
    @(44688, 0, 327144, 6344, 12, 0),
    @(44689, 0, 327157, 6344, 13, 0),
    @(44690, 0, 327167, 6345, 10, 1),
    @(44691, 0, 327179, 6345, 12, 0),
    @(44692, 0, 327198, 6345, 19, 0),
    @(44693, 0, 327218, 6345, 20, 0),
    @(44694, 0, 327238, 6345, 20, 0)
)

$startRow = 36
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 1).NumberFormat = "yyyy\-mm\-dd;@"
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
}

# Update the view to match the final saved state (active selection cell)
$ws.Range("D43").Select()
